$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

# A3, B3, C3, AN3 and AO3 carry a "quote-prefix" (force-text) cell style in
# the source file. A plain .Value assignment of a non-ambiguous alphanumeric
# string causes the engine to re-resolve the cell to the equivalent style
# without the quote-prefix flag, which would shift the style index. Typing
# the value with a leading apostrophe (exactly like typing '... in Excel's
# UI) keeps the text forced and preserves the original style.
$ws.Range("A3").Value = "'JSSO1000251"
$ws.Range("B3").Value = "'JSSO1000251"
$ws.Range("C3").Value = "'JSSO1000251"

$ws.Range("AJ3").Value = "JSCN1000251"
$ws.Range("AL3").Value = "SLJSSO1000251"

$ws.Range("AN3").Value = "'MBLJSSO1000251"
$ws.Range("AO3").Value = "'HBLJSSO1000251"
